$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Create" (sheet1): update row 2 values, selection, and tab state
# ---------------------------------------------------------------------------
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("A2").Value = "ActualSMS"
$wsCreate.Range("B2").Value = "Expected"
$wsCreate.Range("D2").Value = "Chat"
$wsCreate.Range("E2").Value = "'404"
$wsCreate.Range("F2").Value = "'404"
$wsCreate.Range("G2").Value = "Actual Test"
$wsCreate.Range("H2").Value = "'404"

# ---------------------------------------------------------------------------
# Sheet "Edit" (sheet2): update row 2 values and selection
# ---------------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("A2").Value = "ActualSMS"
$wsEdit.Range("B2").Value = "Expected"
$wsEdit.Range("D2").Value = "Chat"
$wsEdit.Range("E2").Value = "'404"
$wsEdit.Range("F2").Value = "'404"
$wsEdit.Range("G2").Value = "Actual Test"
$wsEdit.Range("H2").Value = "'404"
$wsEdit.Range("J2").Value = "Eldorado"

# ---------------------------------------------------------------------------
# Sheet "Delete" (sheet3): update row 2 values and selection
# ---------------------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("A2").Value = "ActualSMS"
$wsDelete.Range("B2").Value = "Expected"
$wsDelete.Range("D2").Value = "Chat"
$wsDelete.Range("E2").Value = "'404"
$wsDelete.Range("F2").Value = "'404"
$wsDelete.Range("G2").Value = "Eldorado"
$wsDelete.Range("H2").Value = "'404"

# ---------------------------------------------------------------------------
# Selections per sheet (set while each sheet is the active one) and final
# active tab. "Create" loses tabSelected, "Delete" gains it (activeTab=2).
# ---------------------------------------------------------------------------
$wsCreate.Activate()
$wsCreate.Range("E2").Select()

$wsEdit.Activate()
$wsEdit.Range("J2").Select()

$wsDelete.Activate()
$wsDelete.Range("G2").Select()
